# Apply cryptos list update (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.601.91"
$ws.Range("E2").Value = "  +2.95%  "
$ws.Range("D3").Value = "2.990.76"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'567.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").Value = "'139.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.97%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").Value = "2.982.99"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").Value = "'0.133"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.47%  "
$ws.Range("E11").Value = "  +11.90%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "'0.0000230"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.15%  "
$ws.Range("D14").Value = "'33.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.88%  "
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "3.481.23"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "'7.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "2.984.82"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").Value = "59.578.88"
$ws.Range("E19").Value = "  +3.01%  "
$ws.Range("D20").Value = "'438.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.33%  "
$ws.Range("D21").Value = "'13.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("D22").Value = "'0.724"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.65%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "'13.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").Value = "'80.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'2.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.81%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("E30").Value = "  +2.86%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'25.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.70%  "
$ws.Range("E33").Value = "  +9.39%  "
$ws.Range("D34").Value = "0.0₃0779"
$ws.Range("E34").Value = "  +11.10%  "
$ws.Range("D35").Value = "'5.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.95%  "
$ws.Range("D36").Value = "'0.983"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.93%  "
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").Value = "'48.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").Value = "'8.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("D41").Value = "'399.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.44%  "
$ws.Range("D42").Value = "'0.0353"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").Value = "2.737.21"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("D45").Value = "'0.251"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.77%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'34.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +20.23%  "
$ws.Range("D48").Value = "'122.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").Value = "'2.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.80%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.110"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.82%  "
$ws.Range("D51").Value = "'23.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.87%  "
